# Restore/update the rule row "R30" (row 10) Integer-min value (column C)
# from 18 to 1 on the "Rules" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
